$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# 1. Title: "ED: " + "Tema" + " 2" -> "ED: Tema 2"
Replace-Text "ED: Tema 2" "ED: Tema 2"

# 2. "El Search Engine Optimization (SEO)" merge
Replace-Text "El Search Engine Optimization (SEO)" "El Search Engine Optimization (SEO)"

# 3. "White-hat" merge (rest of sentence stays a separate run)
Replace-Text "White-hat" "White-hat"

# 4. "Black-hat SEO: utilizar técnicas " merge
Replace-Text "Black-hat SEO: utilizar técnicas " "Black-hat SEO: utilizar técnicas "

# 5. "(Wordpress)" merge
Replace-Text "(Wordpress)" "(Wordpress)"

# 6. "CTR (Click-through rate): porcentaje de usuarios que hacen clic en un enlace." merge
Replace-Text "CTR (Click-through rate): porcentaje de usuarios que hacen clic en un enlace." "CTR (Click-through rate): porcentaje de usuarios que hacen clic en un enlace."

# 7. "GYM: Google, Yahoo y Microsoft. Propietarios de los 3 motores de búsqueda más grandes." merge
Replace-Text "GYM: Google, Yahoo y Microsoft. Propietarios de los 3 motores de búsqueda más grandes." "GYM: Google, Yahoo y Microsoft. Propietarios de los 3 motores de búsqueda más grandes."

# 8. "Link juice: calidad del enlace." merge
Replace-Text "Link juice: calidad del enlace." "Link juice: calidad del enlace."

# 9. "Long tail" merge
Replace-Text "Long tail" "Long tail"

# 10. "SERP (Search Engine Results Page): Las páginas web más relevantes de una consulta." merge
Replace-Text "SERP (Search Engine Results Page): Las páginas web más relevantes de una consulta." "SERP (Search Engine Results Page): Las páginas web más relevantes de una consulta."

# 11. "Stuffing: utilizar la misma palabra clave con demasiada frecuencia" merge
Replace-Text "Stuffing: utilizar la misma palabra clave con demasiada frecuencia" "Stuffing: utilizar la misma palabra clave con demasiada frecuencia"

# 12. Google founders paragraph: text correction "hacer accesible" -> "hacerla accesible" + trailing period
Replace-Text "hacer accesible para todo el mundo" "hacerla accesible para todo el mundo."

# 13. "llamado Webmasters" merge
Replace-Text ". A través de este programa y otro lanzado posteriormente llamado Webmasters" ". A través de este programa y otro lanzado posteriormente llamado Webmasters"

# 14. "Consultas informativas (Know)" merge
Replace-Text "Consultas informativas (Know)" "Consultas informativas (Know)"

# 15. "Consultas de navegación (Go)" merge
Replace-Text "Consultas de navegación (Go)" "Consultas de navegación (Go)"

Write-Output "done"
